$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing municipio/casos/obitos data shifts down.
$ws.Rows.Item(1).Insert()

# Header text
$ws.Range("A1").Value = "MUNICIPIO"
$ws.Range("B1").Value = "CASOS"
$ws.Range("C1").Value = "ÓBITOS"

# Build the header look (bold, thin box border, centered/top aligned) on A1 once,
# then propagate the exact same style to B1:C1 via copy so every header cell
# shares one style record instead of three.
$hdr = $ws.Range("A1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Copy($ws.Range("B1"))
$hdr.Copy($ws.Range("C1"))

# Restore header text after the format copy (Copy also duplicates the source value).
$ws.Range("A1").Value = "MUNICIPIO"
$ws.Range("B1").Value = "CASOS"
$ws.Range("C1").Value = "ÓBITOS"

# Append two new rows at the end of the data (rows 169 and 170)
$ws.Range("A169").Value = "outros estados"
$ws.Range("B169").Value = 43

$ws.Range("A170").Value = "outros paises"
$ws.Range("B170").Value = 42
